# Auto-generated: apply updated leve profit figures from the scheduled runner sync
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(94, 8).Value = 2750.5
$ws.Cells.Item(94, 9).Value = 2667.2222
$ws.Cells.Item(94, 10).Value = 3500
$ws.Cells.Item(94, 11).Value = 2667.2222
$ws.Cells.Item(94, 12).Value = 3500
$ws.Cells.Item(94, 13).Value = -2216.2222
$ws.Cells.Item(94, 14).Value = -4402
$ws.Cells.Item(103, 8).Value = 416
$ws.Cells.Item(103, 9).Value = 374.5
$ws.Cells.Item(103, 10).Value = 499
$ws.Cells.Item(103, 11).Value = 1123.5
$ws.Cells.Item(103, 12).Value = 1497
$ws.Cells.Item(103, 13).Value = -537.5
$ws.Cells.Item(103, 14).Value = -2669
$ws.Cells.Item(113, 8).Value = 2102.9375
$ws.Cells.Item(113, 9).Value = 1800
$ws.Cells.Item(113, 10).Value = 3011.75
$ws.Cells.Item(113, 11).Value = 1800
$ws.Cells.Item(113, 12).Value = 3011.75
$ws.Cells.Item(113, 13).Value = 1454
$ws.Cells.Item(113, 14).Value = -9519.75
$ws.Cells.Item(116, 8).Value = 2144.7368
$ws.Cells.Item(116, 9).Value = 864.2857
$ws.Cells.Item(116, 10).Value = 2891.6667
$ws.Cells.Item(116, 11).Value = 864.2857
$ws.Cells.Item(116, 12).Value = 2891.6667
$ws.Cells.Item(116, 13).Value = 2577.7143
$ws.Cells.Item(116, 14).Value = -9775.6667
$ws.Cells.Item(132, 8).Value = 6212706.5
$ws.Cells.Item(132, 9).Value = 8929959
$ws.Cells.Item(132, 10).Value = 1844.5714
$ws.Cells.Item(132, 11).Value = 26789877
$ws.Cells.Item(132, 12).Value = 5533.7142
$ws.Cells.Item(132, 13).Value = -26787347
$ws.Cells.Item(132, 14).Value = -10593.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2562.2034
$ws.Cells.Item(32, 9).Value = 2281
$ws.Cells.Item(32, 10).Value = 4651.143
$ws.Cells.Item(32, 11).Value = 2281
$ws.Cells.Item(32, 12).Value = 4651.143
$ws.Cells.Item(32, 13).Value = -1994
$ws.Cells.Item(32, 14).Value = -5225.143
$ws.Cells.Item(110, 8).Value = 1994.9445
$ws.Cells.Item(110, 9).Value = 1821.2307
$ws.Cells.Item(110, 10).Value = 2446.6
$ws.Cells.Item(110, 11).Value = 1821.2307
$ws.Cells.Item(110, 12).Value = 2446.6
$ws.Cells.Item(110, 13).Value = 223.7692999999999
$ws.Cells.Item(110, 14).Value = -6536.6
$ws.Cells.Item(132, 8).Value = 13835.777
$ws.Cells.Item(132, 9).Value = 20904.8
$ws.Cells.Item(132, 10).Value = 4999.5
$ws.Cells.Item(132, 11).Value = 62714.39999999999
$ws.Cells.Item(132, 12).Value = 14998.5
$ws.Cells.Item(132, 13).Value = -60184.39999999999
$ws.Cells.Item(132, 14).Value = -20058.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1671.6875
$ws.Cells.Item(20, 9).Value = 1604.1052
$ws.Cells.Item(20, 10).Value = 1770.4615
$ws.Cells.Item(20, 11).Value = 1604.1052
$ws.Cells.Item(20, 12).Value = 1770.4615
$ws.Cells.Item(20, 13).Value = -1357.1052
$ws.Cells.Item(20, 14).Value = -2264.4615
$ws.Cells.Item(62, 8).Value = 33333
$ws.Cells.Item(62, 9).Value = 33333
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 33333
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -32647
$ws.Cells.Item(65, 8).Value = 33333
$ws.Cells.Item(65, 9).Value = 33333
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 99999
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -96567

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1011.1177
$ws.Cells.Item(16, 9).Value = 806.9231
$ws.Cells.Item(16, 10).Value = 1674.75
$ws.Cells.Item(16, 11).Value = 806.9231
$ws.Cells.Item(16, 12).Value = 1674.75
$ws.Cells.Item(16, 13).Value = -519.9231
$ws.Cells.Item(16, 14).Value = -2248.75
$ws.Cells.Item(113, 8).Value = 1011.1177
$ws.Cells.Item(113, 9).Value = 806.9231
$ws.Cells.Item(113, 10).Value = 1674.75
$ws.Cells.Item(113, 11).Value = 806.9231
$ws.Cells.Item(113, 12).Value = 1674.75
$ws.Cells.Item(113, 13).Value = 1363.0769
$ws.Cells.Item(113, 14).Value = -6014.75
$ws.Cells.Item(132, 8).Value = 2501.1428
$ws.Cells.Item(132, 9).Value = 2158.4783
$ws.Cells.Item(132, 11).Value = 6475.4349
$ws.Cells.Item(132, 13).Value = -3945.4349
$ws.Cells.Item(134, 8).Value = 1072.0952
$ws.Cells.Item(134, 9).Value = 1042.1052
$ws.Cells.Item(134, 11).Value = 3126.3156
$ws.Cells.Item(134, 13).Value = -591.3155999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 250
$ws.Cells.Item(7, 9).Value = 250
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 750
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -638
$ws.Cells.Item(7, 14).ClearContents()
$ws.Cells.Item(92, 8).Value = 234.09677
$ws.Cells.Item(92, 9).Value = 121.583336
$ws.Cells.Item(92, 10).Value = 305.1579
$ws.Cells.Item(92, 11).Value = 364.750008
$ws.Cells.Item(92, 12).Value = 915.4737
$ws.Cells.Item(92, 13).Value = 883.249992
$ws.Cells.Item(92, 14).Value = -3411.4737

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1876.3
$ws.Cells.Item(102, 9).Value = 1835.4286
$ws.Cells.Item(102, 10).Value = 1971.6666
$ws.Cells.Item(102, 11).Value = 1835.4286
$ws.Cells.Item(102, 12).Value = 1971.6666
$ws.Cells.Item(102, 13).Value = -213.4286
$ws.Cells.Item(102, 14).Value = -5215.6666
$ws.Cells.Item(126, 8).Value = 2386131.8
$ws.Cells.Item(126, 9).Value = 7002.4
$ws.Cells.Item(126, 10).Value = 3707870.2
$ws.Cells.Item(126, 11).Value = 21007.2
$ws.Cells.Item(126, 12).Value = 11123610.6
$ws.Cells.Item(126, 13).Value = -18537.2
$ws.Cells.Item(126, 14).Value = -11128550.6
$ws.Cells.Item(132, 8).Value = 135402.47
$ws.Cells.Item(132, 9).Value = 168385.33
$ws.Cells.Item(132, 10).Value = 3471
$ws.Cells.Item(132, 11).Value = 505155.99
$ws.Cells.Item(132, 12).Value = 10413
$ws.Cells.Item(132, 13).Value = -502625.99
$ws.Cells.Item(132, 14).Value = -15473

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3099
$ws.Cells.Item(61, 9).Value = 2166.6667
$ws.Cells.Item(61, 11).Value = 2166.6667
$ws.Cells.Item(61, 13).Value = -1964.6667
$ws.Cells.Item(113, 8).Value = 3099
$ws.Cells.Item(113, 9).Value = 2166.6667
$ws.Cells.Item(113, 11).Value = 2166.6667
$ws.Cells.Item(113, 13).Value = 3.333299999999781
$ws.Cells.Item(132, 8).Value = 16093.833
$ws.Cells.Item(132, 9).Value = 30318.182
$ws.Cells.Item(132, 11).Value = 90954.546
$ws.Cells.Item(132, 13).Value = -88424.546
$ws.Cells.Item(138, 8).Value = 20390
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1160.8572
$ws.Cells.Item(122, 10).Value = 1309.5
$ws.Cells.Item(122, 12).Value = 3928.5
$ws.Cells.Item(122, 14).Value = -8828.5
$ws.Cells.Item(126, 8).Value = 1045.091
$ws.Cells.Item(126, 9).Value = 1299.4286
$ws.Cells.Item(126, 10).Value = 600
$ws.Cells.Item(126, 11).Value = 3898.2858
$ws.Cells.Item(126, 12).Value = 1800
$ws.Cells.Item(126, 13).Value = -1428.2858
$ws.Cells.Item(126, 14).Value = -6740
$ws.Cells.Item(132, 8).Value = 2470.0715
$ws.Cells.Item(132, 9).Value = 1683.2858
$ws.Cells.Item(132, 11).Value = 5049.857400000001
$ws.Cells.Item(132, 13).Value = -2519.857400000001

